$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.510.55"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").Value = "1.915.64"
$ws.Range("E3").Value = "  -0.38%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4794"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2891"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.41%  "

$ws.Range("E9").Value = "  -0.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "110.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.80%  "

$ws.Range("D12").Value = "1.910.57"
$ws.Range("E12").Value = "  +1.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07570"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.246"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6672"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "301.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.02%  "

$ws.Range("D17").Value = "30.493.64"
$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9994"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007566"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.72%  "

$ws.Range("D21").Value = "2.161.44"
$ws.Range("E21").Value = "  +1.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.490"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.418"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.43%  "

$ws.Range("E25").Value = "  +0.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.102"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1075"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.393"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.04%  "

$ws.Range("E31").Value = "  -0.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.017"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04983"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7363"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.47%  "

$ws.Range("E35").Value = "  -1.59%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02049"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9983"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("E38").Value = "  -0.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.676"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.37%  "

$ws.Range("E41").Value = "  -2.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4431"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8652"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.902"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.07%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9996"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "50.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.278"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.309"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1232"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2522"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.40%  "
